$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $r = $ws.Range($addr)
    $r.NumberFormat = "@"
    $r.Value = $val
    $r.ClearFormats()
}

# Row 2 - Bitcoin
Set-TextValue "D2" "25.928.73"
$ws.Range("E2").Value = "  -1.09%  "

# Row 3 - Ethereum
Set-TextValue "D3" "1.636.49"
$ws.Range("E3").Value = "  -0.49%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.25%  "

# Row 5 - BNB
Set-TextValue "D5" "215.36"
$ws.Range("E5").Value = "  -0.70%  "

# Row 6 - XRP
Set-TextValue "D6" "0.507"
$ws.Range("E6").Value = "  +0.19%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  +0.17%  "

# Row 8 - Cardano
$ws.Range("E8").Value = "  -0.82%  "

# Row 9 - Dogecoin
Set-TextValue "D9" "0.0636"
$ws.Range("E9").Value = "  -0.11%  "

# Row 10 - Solana
$ws.Range("E10").Value = "  -1.44%  "

# Row 11 - TRON
$ws.Range("E11").Value = "  +0.17%  "

# Row 12 - Polkadot
$ws.Range("E12").Value = "  -0.08%  "

# Row 13 - WrappedliquidstakedEther2.0
Set-TextValue "D13" "1.863.05"
$ws.Range("E13").Value = "  -0.49%  "

# Row 14 - WrappedEther
Set-TextValue "D14" "1.609.62"
$ws.Range("E14").Value = "  -2.75%  "

# Row 15 - Polygon
$ws.Range("E15").Value = "  -0.54%  "

# Row 16 - ShibaInu
$ws.Range("E16").Value = "  -0.22%  "

# Row 17 - Litecoin
$ws.Range("E17").Value = "  -0.64%  "

# Row 18 - WrappedBTC
Set-TextValue "D18" "25.889.13"
$ws.Range("E18").Value = "  -1.24%  "

# Row 19 - Dai
$ws.Range("E19").Value = "  +0.22%  "

# Row 20 - BitcoinCash
Set-TextValue "D20" "192.79"
$ws.Range("E20").Value = "  -1.27%  "

# Row 21 - Uniswap
$ws.Range("E21").Value = "  -1.79%  "

# Row 22 - Avalanche
$ws.Range("E22").Value = "  -1.38%  "

# Row 23 - Chainlink
$ws.Range("E23").Value = "  -0.40%  "

# Row 24 - Toncoin
Set-TextValue "D24" "1.80"
$ws.Range("E24").Value = "  +0.95%  "

# Row 25 / Row 26 - Stellar and Monero swap places (content swap, same rows)
$ws.Range("B25").Value = "Monero"
$ws.Range("C25").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue "D25" "144.15"
$ws.Range("E25").Value = "  +1.03%  "

$ws.Range("B26").Value = "Stellar"
$ws.Range("C26").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
Set-TextValue "D26" "0.131"
$ws.Range("E26").Value = "  +4.65%  "

# Row 27 - BinanceUSD
$ws.Range("E27").Value = "  +0.27%  "

# Row 28 - Cosmos
Set-TextValue "D28" "6.90"
$ws.Range("E28").Value = "  -0.58%  "

# Row 29 - EthereumClassic
$ws.Range("E29").Value = "  -0.57%  "

# Row 30 - PancakeSwap
$ws.Range("E30").Value = "  -0.41%  "

# Row 31 - Hedera
$ws.Range("E31").Value = "  -0.37%  "

# Row 32 - InternetComputer(DFINITY)
$ws.Range("E32").Value = "  -2.05%  "

# Row 33 - Filecoin
$ws.Range("E33").Value = "  -0.15%  "

# Row 34 - LidoDAOToken
$ws.Range("E34").Value = "  -3.86%  "

# Row 35 - HuobiToken
Set-TextValue "D35" "2.46"
$ws.Range("E35").Value = "  +1.80%  "

# Row 36 - ARBITRUM
$ws.Range("E36").Value = "  -1.22%  "

# Row 37 - Maker (only price changes, volume stays -0.16%)
Set-TextValue "D37" "1.132.23"

# Row 38 - ImmutableX
$ws.Range("E38").Value = "  -1.57%  "

# Row 39 - MXToken
$ws.Range("E39").Value = "  -1.80%  "

# Row 40 - VeChain
$ws.Range("E40").Value = "  -0.43%  "

# Row 41 - FraxShare
Set-TextValue "D41" "5.49"
$ws.Range("E41").Value = "  +0.32%  "

# Row 42 - Quant
Set-TextValue "D42" "99.45"
$ws.Range("E42").Value = "  -0.95%  "

# Row 43 - TrustWalletToken
Set-TextValue "D43" "0.793"
$ws.Range("E43").Value = "  -0.45%  "

# Row 44 - RocketPoolETH
Set-TextValue "D44" "1.772.56"
$ws.Range("E44").Value = "  -0.51%  "

# Row 45 - BabyDogeCoin
Set-TextValue "D45" "0.0₆0114"
$ws.Range("E45").Value = "  +2.74%  "

# Row 46 - Aave
Set-TextValue "D46" "56.61"
$ws.Range("E46").Value = "  -0.64%  "

# Row 47 - Cronos
Set-TextValue "D47" "0.0530"
$ws.Range("E47").Value = "  +2.48%  "

# Row 48 - RenderToken
$ws.Range("E48").Value = "  -0.92%  "

# Row 49 - EnergySwap
Set-TextValue "D49" "7.68"
$ws.Range("E49").Value = "  -0.41%  "

# Row 50 - Mantle
$ws.Range("E50").Value = "  -0.77%  "

# Row 51 - Algorand
Set-TextValue "D51" "0.0958"
$ws.Range("E51").Value = "  -1.13%  "
